$wb = $excel.ActiveWorkbook

# Rename the first sheet (this also updates the _xlnm._FilterDatabase
# defined name's sheet reference automatically)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "LH_TC_IDCONSTRAINSTS"

# Make sheet1 the active/selected sheet: this moves tabSelected="1" from
# sheet2's sheetView to sheet1's sheetView, and the workbook-level
# activeTab goes back to sheet1 (index 0), i.e. the activeTab attribute
# is effectively removed from workbookView.
$ws1.Activate()
$ws1.Select()
